# Update "想去人数" (interested-count) values on the "展览" and "全部类型"
# sheets to match the newly generated output (commit 456a3b4).
#
# Both sheets list the same events, but interleaved with extra rows on
# "全部类型" (it also contains 演出/本地生活 entries), so the target rows
# are looked up by the event name in column C rather than by a fixed
# row number.

$wb = $excel.ActiveWorkbook

# Map of event name -> new "想去人数" value
$updates = @{
    "常熟·ACG动漫游戏嘉年华" = 559
    "常熟·CDW.动漫展05" = 1591
    "苏州·星部落动漫嘉年华" = 5212
    "苏州·国乙主场·次元燃歌·偶像时刻（免费展）" = 204
    "苏州·授渔仲夏动漫节2.0" = 780
    "张家港·META萌圆饿了" = 69
    "昆山·2024首届华盟次元动漫嘉年华（免费活动）" = 377
    "苏州·星河璀璨，爱恋永恒——七夕CP漫游记暨坏孩纸物语NO.48（免费的漫展）" = 528
    "苏州·甜蜜元气偶像日" = 36
    "苏州·第三届.OCG.Summer Carnival-国潮动漫游戏嘉年华" = 6598
    "苏州·艾卡动漫游戏嘉年华（免票展）" = 34
    "昆山·第七届·xcy新次元动漫嘉年华-狂欢盛典" = 143
    "昆山·创世次元动漫游戏嘉年华" = 1039
    "苏州·ICAN summer World动漫品牌夏游节" = 15805
    "苏州·星梦X动漫游戏展（免费展）" = 158
    "苏州·Good jump ACG中秋嘉年华动漫国潮文化节" = 11188
    "苏州·理想乡动漫游戏展-两馆全开+三馆间通道" = 4389
    "苏州·第二届百合Only同人展交流" = 28
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # find last used row
    $lastRow = $ws.UsedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        # NOTE: reading via ".Value" is unreliable in this runtime, use
        # ".Text" instead (writing via ".Value" works fine).
        $name = $ws.Cells.Item($r, 3).Text
        if ($null -ne $name -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
